# Auto-generated Excel COM-interop script
# Updates column F (想去人数 / 'interested count') values across sheets
# as per the diff between before.xlsx and the target output.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 1513
$ws.Cells.Item(3, 6).Value = 851
$ws.Cells.Item(5, 6).Value = 888
$ws.Cells.Item(6, 6).Value = 502
$ws.Cells.Item(7, 6).Value = 7549
$ws.Cells.Item(11, 6).Value = 5476
$ws.Cells.Item(15, 6).Value = 7549
$ws.Cells.Item(16, 6).Value = 8893
$ws.Cells.Item(17, 6).Value = 1141
$ws.Cells.Item(18, 6).Value = 887
$ws.Cells.Item(19, 6).Value = 4421
$ws.Cells.Item(20, 6).Value = 665
$ws.Cells.Item(21, 6).Value = 214
$ws.Cells.Item(25, 6).Value = 1191
$ws.Cells.Item(26, 6).Value = 105
$ws.Cells.Item(27, 6).Value = 1649
$ws.Cells.Item(28, 6).Value = 705
$ws.Cells.Item(29, 6).Value = 902
$ws.Cells.Item(30, 6).Value = 4
$ws.Cells.Item(31, 6).Value = 1862
$ws.Cells.Item(32, 6).Value = 332
$ws.Cells.Item(33, 6).Value = 2271
$ws.Cells.Item(35, 6).Value = 107
$ws.Cells.Item(36, 6).Value = 1446
$ws.Cells.Item(39, 6).Value = 790
$ws.Cells.Item(40, 6).Value = 406
$ws.Cells.Item(41, 6).Value = 4068
$ws.Cells.Item(42, 6).Value = 191
$ws.Cells.Item(47, 6).Value = 858
$ws.Cells.Item(49, 6).Value = 4075

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(12, 6).Value = 27

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 5182

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(4, 6).Value = 1513
$ws.Cells.Item(5, 6).Value = 851
$ws.Cells.Item(6, 6).Value = 888
$ws.Cells.Item(7, 6).Value = 502
$ws.Cells.Item(11, 6).Value = 5476
$ws.Cells.Item(13, 6).Value = 7549
$ws.Cells.Item(14, 6).Value = 27
$ws.Cells.Item(15, 6).Value = 1141
$ws.Cells.Item(16, 6).Value = 887
$ws.Cells.Item(17, 6).Value = 4421
$ws.Cells.Item(18, 6).Value = 665
$ws.Cells.Item(19, 6).Value = 214
$ws.Cells.Item(25, 6).Value = 1191
$ws.Cells.Item(26, 6).Value = 105
$ws.Cells.Item(27, 6).Value = 1649
$ws.Cells.Item(28, 6).Value = 1862
$ws.Cells.Item(29, 6).Value = 332
$ws.Cells.Item(30, 6).Value = 2271
$ws.Cells.Item(37, 6).Value = 790
$ws.Cells.Item(40, 6).Value = 406
$ws.Cells.Item(41, 6).Value = 4068
$ws.Cells.Item(43, 6).Value = 191
$ws.Cells.Item(47, 6).Value = 858
$ws.Cells.Item(49, 6).Value = 4075
